# Apply the "7.2.1" -> "7.2.1.1" indicator renumbering plus the few data
# corrections shown in the diff (sheet1 of the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cells: renumber "7.2.1" -> "7.2.1.1" --------------------------
$ws.Range("A1").Value = " 7.2.1.1 Энергия керектөөлөрүнүн жалпы көлөмүндөгү энергиянын жаңыланган булактарынын  үлүшү"
$ws.Range("B1").Value = " 7.2.1.1 Доля возобновляемых источников энергии в общем объеме энергопотребления"
$ws.Range("C1").Value = "7.2.1.1 Renewable energy share in the total energy consumption"

# --- Data corrections in rows 5 & 6 ---------------------------------------
$ws.Range("Q5").Value = 36.700000000000003
$ws.Range("P6").Value = 13859.2
$ws.Range("Q6").Value = 13979.2

# --- Selection moved from P9 to P7 ----------------------------------------
$ws.Range("P7").Select()
